$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Metadata sheet ("Metadata") updates for RAD IMR 1.1.0 build ---

# Version: 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Experimental: (empty) -> "false"
# Plain assignment would auto-type "false" as a Boolean; round-trip it
# through a text formula + paste-values so it lands as literal text,
# matching the original cell's string-based style.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Date: refreshed build timestamp
$ws.Range("B8").Value = "2024-06-20T08:51:57-05:00"

# Contact rows: three ContactDetail entries now resolve to real display text
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/radiology/)"
$ws.Range("B11").Value = "null (radiology@ihe.net)"
$ws.Range("B12").Value = "IHE Radiology Technical Committee (radiology@ihe.net)"

# Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"
